$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column C header, copying B1's formatting (style 2)
$ws.Range("C1").Value = "comments"
$ws.Range("B1").Copy()
$ws.Range("C1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Set column C width (closest achievable to target 42.796875 given engine rounding)
$ws.Columns.Item(3).ColumnWidth = 42.14

# Comment for row 7 (Water industry or related sectors experience)
$ws.Range("C7").Value = "Water industry experience is preferred`nUtilitary experience is mandatory"
$ws.Range("C7").WrapText = $true

# Comment for row 8 (Key personnel proposed to provide the service)
$ws.Range("C8").Value = "The backgrounds of key personnel must be provided`nAn indication of the availability of key personnel must be provided`nPersonnel must be appropriately qualified."
$ws.Range("C8").WrapText = $true
$ws.Rows.Item(8).RowHeight = 69

# Update selection to reflect final active cell in the original edit session
$ws.Range("C18").Select()

